$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three changed data values in row 2 (E2, G2, H2)
$ws.Range("E2").Value = 6
$ws.Range("G2").Value = -3
$ws.Range("H2").Value = 13

# Move/record the active selection to E2, matching the saved view state
$ws.Range("E2").Select()
